$wb = $excel.ActiveWorkbook

# --- "About" sheet: add the two new rows (35 and 36) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Cells.Item(35,1).Value = 7.8285
$wsAbout.Cells.Item(35,2).Value = "USD converted to HKD"
$wsAbout.Cells.Item(36,1).Formula = '=A26*A35'
$wsAbout.Cells.Item(36,2).Value = "this number was used"

# --- Dependent OCCF sheets now point at the new row 36 instead of row 26 ---
$wsLOCU = $wb.Worksheets.Item("OCCF-DpLOCU")
$wsLOCU.Range("B2").Formula = '=10^9*About!$A$36'
$wsLOCU.Range("B2").NumberFormat = "0"

$wsMOCU = $wb.Worksheets.Item("OCCF-DpMOCU")
$wsMOCU.Range("B2").Formula = '=10^6*About!$A$36'

$wsSOCU = $wb.Worksheets.Item("OCCF-DpSOCU")
$wsSOCU.Range("B2").Formula = '=1*About!A36'

# --- Update selections on every sheet ---
$null = $wsAbout.Range("A36").Select()

$null = $wsLOCU.Range("B2").Select()
$null = $wsMOCU.Range("B2").Select()
$null = $wsSOCU.Range("B2").Select()

# --- Active sheet/tab moves from "About" to "OCCF-DpLOCU" ---
$null = $wsLOCU.Activate()
$null = $wsLOCU.Range("B2").Select()
